$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8534.616716272223
$ws.Range("C2").Value = 25690.62581580083
$ws.Range("D2").Value = 60654.61806333878
$ws.Range("E2").Value = 100531.5535416491

$ws.Range("B3").Value = 89885.81977983763
$ws.Range("C3").Value = 232306.1762703886
$ws.Range("D3").Value = 339991.8627747862
$ws.Range("E3").Value = 413840.7570057246

$ws.Range("B4").Value = 10090.6497544036
$ws.Range("C4").Value = 26190.64920742056
$ws.Range("D4").Value = 49951.67028788106
$ws.Range("E4").Value = 70609.47614285229

$ws.Range("B6").Value = 52984.19775889564
$ws.Range("C6").Value = 93751.81611425432
$ws.Range("D6").Value = 103884.7425147522
$ws.Range("E6").Value = 93137.31223794348

$ws.Range("B7").Value = 5770.913820418582
$ws.Range("C7").Value = 15644.51171895376
$ws.Range("D7").Value = 20529.88086988642
$ws.Range("E7").Value = 24365.92197415622

$ws.Range("B9").Value = 419705.8813295108
$ws.Range("C9").Value = 918505.808206992
$ws.Range("D9").Value = 1487443.043026546
$ws.Range("E9").Value = 1999620.638853451

$ws.Range("B12").Value = 561717.0829547446
$ws.Range("C12").Value = 863172.7447987135
$ws.Range("D12").Value = 870283.3226981713
$ws.Range("E12").Value = 689812.7500598714
